$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "QuixBugs(Java)"

[void]$ws.Range("I20").Select()
